# Update countries & provincias Spain
# - Refresh the "Casos totales / Nuevos casos / Casos activos / Recuperados /
#   Casos criticos / Muertes hoy / Muertes" figures for a handful of
#   countries (the underlying "Pais" table is kept sorted by "Casos
#   totales" descending, so Costa Rica overtaking Etiopia re-orders those
#   two rows).
# - Bump the "Datos actualizados..." timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowVals($ws, $Row, $map) {
    foreach ($col in $map.Keys) {
        $ws.Cells.Item($Row, $col).Value = $map[$col]
    }
}

# Column numbers: A=1 B=2 C=3 D=4 E=5 F=6 G=7 H=8

# Estados Unidos
Set-RowVals $ws 4 @{ 2 = 4955349; 3 = 36929; 4 = 2509377; 5 = 2284773; 7 = 909; 8 = 161199 }

# Alemania
Set-RowVals $ws 21 @{ 2 = 214104; 3 = 1024; 5 = 10159; 7 = 13; 8 = 9245 }

# Israel
Set-RowVals $ws 36 @{ 2 = 77919; 3 = 1721; 4 = 51395; 5 = 25959 }

# Argelia
Set-RowVals $ws 60 @{ 4 = 22802; 5 = 8992 }

# Uzbekistan
Set-RowVals $ws 62 @{ 2 = 27793; 3 = 746; 4 = 18783; 5 = 8839; 7 = 6; 8 = 171 }

# Costa Rica now has more total cases than Etiopia, so it moves above it in
# the sorted table: row 70 becomes Costa Rica (new figures), row 71 becomes
# Etiopia (its figures are unchanged from before, it just shifts down).
Set-RowVals $ws 70 @{ 1 = "Costa Rica"; 2 = 20417; 3 = 580; 4 = 6851; 5 = 13375; 7 = 10; 8 = 191 }
Set-RowVals $ws 71 @{ 1 = "Etiopia"; 2 = 20336; 3 = 459; 4 = 8598; 5 = 11382; 7 = 13; 8 = 356 }

# Costa de Marfil
Set-RowVals $ws 76 @{ 2 = 16349; 3 = 56; 4 = 12191; 5 = 4055 }

# Mauritania
Set-RowVals $ws 97 @{ 2 = 6444; 3 = 26; 4 = 5291; 5 = 996 }

# Grecia
Set-RowVals $ws 103 @{ 2 = 4974; 3 = 119; 5 = 3390 }

# Ruanda
Set-RowVals $ws 129 @{ 2 = 2104; 3 = 5; 4 = 1237; 5 = 862 }

# Angola
Set-RowVals $ws 139 @{ 2 = 1395; 3 = 51; 4 = 506; 5 = 827; 7 = 3; 8 = 62 }

# Monaco
Set-RowVals $ws 187 @{ 2 = 125; 3 = 2; 5 = 16 }

# Santa Lucia / Timor Oriental are tied on every figure, they simply swap
# their display order in the table.
Set-RowVals $ws 202 @{ 1 = "Santa Lucia" }
Set-RowVals $ws 203 @{ 1 = "Timor Oriental" }

# Refresh the "updated at" banner.
$ws.Range("A1").Value = "Datos actualizados a 5 de Agosto de 2020 a las 22:12"
